# Automatische test-sync: 2025-06-17 22:01:01
# Append a new "Afmelding nieuwsbrief" log entry to the Logs sheet,
# extend the conditional-formatting ranges to cover the new row, and
# bump the "Afmelding" tally on the Dashboard sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$newRow = 40

$logs.Cells.Item($newRow, 1).Value = "Afmelding nieuwsbrief"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$logs.Cells.Item($newRow, 4).Value = "Afmelding"
$logs.Cells.Item($newRow, 6).Value = "2025-06-17 22:00:16"
$logs.Cells.Item($newRow, 7).Value = "Nee"

# Extend the existing conditional formatting (Categorie + Beantwoord
# columns) so it keeps covering the whole data range including the row
# that was just appended.
$catFormatting = $logs.Range("D2:D39").FormatConditions
$catFormatting.Item(1).ModifyAppliesToRange($logs.Range("D2:D40"))

$repliedFormatting = $logs.Range("G2:G39").FormatConditions
$repliedFormatting.Item(1).ModifyAppliesToRange($logs.Range("G2:G40"))

# Bump the "Afmelding" tally on the Dashboard sheet (row 4) from 5 to 6.
$dashboard.Range("B4").Value = 6
